# Move flux bounds (Min flux / Max flux) from the "Rate laws" sheet to the
# "Reactions" sheet, and drop the now-unused "Forward/Reverse rate law"
# columns and the stray "biomass" row in "Species types" (it isn't a real
# metabolite / species type entry).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Species types": remove the bogus "biomass" row (row 8) -- biomass is
#    a reaction, not a species type / metabolite.
# ---------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species types")
$wsSpecies.Range("A8:J8").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2. "Reactions": insert two new columns (G:H) to host the flux bounds that
#    used to live on the "Rate laws" sheet, then drop the old "Forward rate
#    law" / "Reverse rate law" columns (now shifted to L:M).
# ---------------------------------------------------------------------
$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Range("G1:H1").EntireColumn.Insert()
$wsReactions.Range("L1:M1").EntireColumn.Delete()

$wsReactions.Range("G1").Value = "Min flux"
$wsReactions.Range("H1").Value = "Max flux"

$wsReactions.Range("G2").Value = 0
$wsReactions.Range("H2").Value = 1

$wsReactions.Range("G3").Value = 1
$wsReactions.Range("H3").Value = 2

$wsReactions.Range("G4").Value = 2
$wsReactions.Range("H4").Value = 3

$wsReactions.Range("G5").Value = 3
$wsReactions.Range("H5").Value = 4

$wsReactions.Range("G6").Value = 4
$wsReactions.Range("H6").Value = 5

# row 7 (ex_specie_2) keeps no flux bounds, matching the old "Rate laws" data

# ---------------------------------------------------------------------
# 3. "Rate laws": the flux bounds moved away, delete the old F:G columns.
# ---------------------------------------------------------------------
$wsRateLaws = $wb.Worksheets.Item("Rate laws")
$wsRateLaws.Range("F1:G1").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 4. Restore reasonable selections / active sheet, mirroring the manual
#    review pass implied by the commit message ("check bounds & biomass
#    reaction").
# ---------------------------------------------------------------------
$wsSpecies.Activate()
$wsSpecies.Rows(8).Select()

$wsReactions.Activate()
$wsReactions.Range("A6").Select()

$wsRateLaws.Activate()
$wsRateLaws.Range("F1:G1").EntireColumn.Select()
